$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The custom date numFmt (used by the A/B "date" columns) now also shows
# hours/minutes instead of just the day/month/year.
$ws.Range("A2:B9").NumberFormat = "dd.mm.yyyy HH:MM"

# Row 7: transaction/posting dates moved forward (raw Excel date serials).
$ws.Range("A7").Value = 43682.71875
$ws.Range("B7").Value = 43683

# Row 8: transaction/posting dates moved forward.
$ws.Range("A8").Value = 43681.71875
$ws.Range("B8").Value = 43682

# Row 9: only the posting date changes; the transaction date is unchanged.
$ws.Range("B9").Value = 43682
